# Replace the static "www.drpaulduenas.com" text in the footer with a
# configurable MERGEFIELD ("=website"), matching the pattern already used
# by the other footer fields (address, phone number, city, emergency
# number), i.e. fldChar begin / instrText / fldChar separate / result
# text / fldChar end - all wrapped in the SAME run formatting that the
# literal text run used to carry.

$d = $word.ActiveDocument

# The footer paragraph lives in the primary ("default") footer of the
# (only) section.
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

# Locate the literal text inside the footer range.
$rng = $ftr.Range.Duplicate
$found = $rng.Find.Execute("www.drpaulduenas.com", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'www.drpaulduenas.com' in the footer"
}

# Capture the run formatting (rPr) so the new runs match exactly what the
# literal-text run used to carry.
$rPr = '<w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

# Remove the old literal text, collapsing the range to the insertion point.
$rng.Text = ""

# Build the replacement runs: fldChar begin -> instrText -> fldChar
# separate -> cached result text -> fldChar end.
$openChevron = [char]0x00AB
$closeChevron = [char]0x00BB

$runBegin    = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="begin"/></w:r>'
$runInstr    = '<w:r>' + $rPr + '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>'
$runSep      = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="separate"/></w:r>'
$runResult   = '<w:r>' + $rPr + '<w:t>' + $openChevron + '=website' + $closeChevron + '</w:t></w:r>'
$runEnd      = '<w:r>' + $rPr + '<w:fldChar w:fldCharType="end"/></w:r>'

# Re-wrap the target paragraph's own identity/properties so InsertXML
# (which replaces the enclosing paragraph when the range sits fully
# inside it) doesn't drop them.
$pOpen = '<w:p w14:paraId="24EA949D" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRDefault="004D2A29" w:rsidP="004D2A29">'
$pPr   = '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>'
$pClose = '</w:p>'

$newParagraph = $pOpen + $pPr + $runBegin + $runInstr + $runSep + $runResult + $runEnd + $pClose

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $newParagraph + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData>' + `
    '</pkg:part>' + `
    '</pkg:package>'

$rng.InsertXML($packageXml)

Write-Output "Replaced website text with MERGEFIELD in footer: $($ftr.Range.Text)"
